$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("all_tools")
$ws5 = $wb.Worksheets.Item("openjml")

# --- Updates to all_tools sheet (rows 13-24) ---
$ws1.Range("G13").Value = 36
$ws1.Range("I13").Value = 0.04969039949999533
$ws1.Range("J13").Value = 0.8509806870320558
$ws1.Range("K13").Value = 0.1010313640565443
$ws1.Range("L13").Value = 0.7812359924917622
$ws1.Range("G14").Value = 36
$ws1.Range("I14").Value = 0.3975231959999626
$ws1.Range("J14").Value = 0.1328549557310538
$ws1.Range("K14").Value = 0.4925278997756535
$ws1.Range("L14").Value = 0.148112135570859
$ws1.Range("G15").Value = 36
$ws1.Range("I15").Value = 0.4773929622481257
$ws1.Range("J15").Value = 0.07328176503287165
$ws1.Range("K15").Value = 0.5890324631280877
$ws1.Range("L15").Value = 0.07317267760755417
$ws1.Range("G16").Value = 36
$ws1.Range("I16").Value = -0.4969039949999532
$ws1.Range("J16").Value = 0.06028917399060209
$ws1.Range("K16").Value = -0.5304146612968577
$ws1.Range("L16").Value = 0.1147392659290222
$ws1.Range("G17").Value = 36
$ws1.Range("I17").Value = -0.6956655929999345
$ws1.Range("J17").Value = 0.008534920414227074
$ws1.Range("K17").Value = -0.8335087534664906
$ws1.Range("L17").Value = 0.002735455303093727
$ws1.Range("G18").Value = 36
$ws1.Range("I18").Value = 0.04969039949999533
$ws1.Range("J18").Value = 0.8509806870320558
$ws1.Range("K18").Value = 0.02525784101413608
$ws1.Range("L18").Value = 0.9447837074747326
$ws1.Range("G19").Value = 36
$ws1.Range("I19").Value = 0.07537783614444091
$ws1.Range("J19").Value = 0.7773295263554205
$ws1.Range("K19").Value = 0.04433577679458724
$ws1.Range("L19").Value = 0.9032059022108545
$ws1.Range("G20").Value = 36
$ws1.Range("I20").Value = 0
$ws1.Range("J20").Value = 1
$ws1.Range("K20").Value = 0.01262892050706804
$ws1.Range("L20").Value = 0.9723786419920799
$ws1.Range("G21").Value = 36
$ws1.Range("I21").Value = -0.4969039949999532
$ws1.Range("J21").Value = 0.06028917399060209
$ws1.Range("K21").Value = -0.631446025353402
$ws1.Range("L21").Value = 0.05021407909522695
$ws1.Range("G22").Value = 36
$ws1.Range("I22").Value = 0.5962847939999439
$ws1.Range("J22").Value = 0.02417054717454525
$ws1.Range("K22").Value = 0.7198484689028782
$ws1.Range("L22").Value = 0.01890477781850608
$ws1.Range("G23").Value = 36
$ws1.Range("I23").Value = 0.5465943944999486
$ws1.Range("J23").Value = 0.03877750439230662
$ws1.Range("K23").Value = 0.6945906278887423
$ws1.Range("L23").Value = 0.0258112209674474
$ws1.Range("G24").Value = 36
$ws1.Range("I24").Value = -0.6956655929999345
$ws1.Range("J24").Value = 0.008534920414227074
$ws1.Range("K24").Value = -0.8208798329594226
$ws1.Range("L24").Value = 0.003605943791192453
# --- Updates to openjml sheet (rows 13-24) ---
$ws5.Range("F13").Value = 10
$ws5.Range("G13").Value = 26
$ws5.Range("I13").Value = 0.04969039949999533
$ws5.Range("J13").Value = 0.8509806870320558
$ws5.Range("K13").Value = 0.01262892050706804
$ws5.Range("L13").Value = 0.9723786419920799
$ws5.Range("F14").Value = 10
$ws5.Range("G14").Value = 26
$ws5.Range("I14").Value = -0.149071198499986
$ws5.Range("J14").Value = 0.5730251193553904
$ws5.Range("K14").Value = -0.1894338076060206
$ws5.Range("L14").Value = 0.6001664342511973
$ws5.Range("F15").Value = 10
$ws5.Range("G15").Value = 26
$ws5.Range("I15").Value = -0.1256297269074015
$ws5.Range("J15").Value = 0.6374017405958849
$ws5.Range("K15").Value = -0.152008377581442
$ws5.Range("L15").Value = 0.6750590889374006
$ws5.Range("F16").Value = 10
$ws5.Range("G16").Value = 26
$ws5.Range("I16").Value = -0.298142396999972
$ws5.Range("J16").Value = 0.2596563563704499
$ws5.Range("K16").Value = -0.3788676152120412
$ws5.Range("L16").Value = 0.2802942824523375
$ws5.Range("F17").Value = 10
$ws5.Range("G17").Value = 26
$ws5.Range("I17").Value = -0.1987615979999813
$ws5.Range("J17").Value = 0.4523703606773608
$ws5.Range("K17").Value = -0.3409808536908371
$ws5.Range("L17").Value = 0.3349456951179903
$ws5.Range("F18").Value = 10
$ws5.Range("G18").Value = 26
$ws5.Range("I18").Value = 0.149071198499986
$ws5.Range("J18").Value = 0.5730251193553904
$ws5.Range("K18").Value = 0.1641759665918845
$ws5.Range("L18").Value = 0.6503895621649565
$ws5.Range("F19").Value = 10
$ws5.Range("G19").Value = 26
$ws5.Range("I19").Value = 0.07537783614444091
$ws5.Range("J19").Value = 0.7773295263554205
$ws5.Range("K19").Value = 0.1076726007868547
$ws5.Range("L19").Value = 0.7671778789420547
$ws5.Range("F20").Value = 10
$ws5.Range("G20").Value = 26
$ws5.Range("I20").Value = -0.09938079899999065
$ws5.Range("J20").Value = 0.7071142312899612
$ws5.Range("K20").Value = -0.1262892050706804
$ws5.Range("L20").Value = 0.7281063840216824
$ws5.Range("F21").Value = 10
$ws5.Range("G21").Value = 26
$ws5.Range("I21").Value = -0.149071198499986
$ws5.Range("J21").Value = 0.5730251193553904
$ws5.Range("K21").Value = -0.2399494896342928
$ws5.Range("L21").Value = 0.5043017190353258
$ws5.Range("F22").Value = 10
$ws5.Range("G22").Value = 26
$ws5.Range("I22").Value = 0.3975231959999626
$ws5.Range("J22").Value = 0.1328549557310538
$ws5.Range("K22").Value = 0.5304146612968577
$ws5.Range("L22").Value = 0.1147392659290222
$ws5.Range("F23").Value = 10
$ws5.Range("G23").Value = 26
$ws5.Range("I23").Value = 0.298142396999972
$ws5.Range("J23").Value = 0.2596563563704499
$ws5.Range("K23").Value = 0.4041254562261773
$ws5.Range("L23").Value = 0.2467547295422347
$ws5.Range("F24").Value = 10
$ws5.Range("G24").Value = 26
$ws5.Range("I24").Value = -0.2484519974999766
$ws5.Range("J24").Value = 0.347558036741169
$ws5.Range("K24").Value = -0.4293832972403134
$ws5.Range("L24").Value = 0.2155824117700313
Write-Host "Applied correlation analysis updates to all_tools and openjml sheets"
